$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "67.813.35"
$ws.Range("E2").Value = "  +0.59%  "

$ws.Range("D3").Value = "2.497.16"
$ws.Range("E3").Value = "  -2.22%  "

$ws.Range("E4").Value = "  +0.01%  "

$ws.Range("D5").Value2 = "'592.30"
$ws.Range("D5").Style = "Normal"
$ws.Range("E5").Value = "  -0.26%  "

$ws.Range("D6").Value2 = "'174.01"
$ws.Range("D6").Style = "Normal"
$ws.Range("E6").Value = "  +0.42%  "

$ws.Range("E7").Value = "  -0.05%  "

$ws.Range("E8").Value = "  -1.10%  "

$ws.Range("D9").Value = "2.495.67"
$ws.Range("E9").Value = "  -2.22%  "

$ws.Range("E10").Value = "  -0.05%  "

$ws.Range("E11").Value = "  +1.73%  "

$ws.Range("D12").Value2 = "'5.11"
$ws.Range("D12").Style = "Normal"
$ws.Range("E12").Value = "  -1.33%  "

$ws.Range("D13").Value2 = "'0.343"
$ws.Range("D13").Style = "Normal"
$ws.Range("E13").Value = "  -2.80%  "

$ws.Range("D14").Value2 = "'26.30"
$ws.Range("D14").Style = "Normal"
$ws.Range("E14").Value = "  -3.10%  "

$ws.Range("D15").Value = "2.951.06"
$ws.Range("E15").Value = "  -2.09%  "

$ws.Range("E16").Value = "  -1.29%  "

$ws.Range("D17").Value = "67.717.85"
$ws.Range("E17").Value = "  +0.71%  "

$ws.Range("D18").Value = "2.493.03"
$ws.Range("E18").Value = "  -2.49%  "

$ws.Range("B19").Value = "Chainlink"
$ws.Range("C19").Value = "https://coinranking.com/coin/VLqpJwogdhHNb+chainlink-link"
$ws.Range("D19").Value2 = "'11.72"
$ws.Range("D19").Style = "Normal"
$ws.Range("E19").Value = "  +2.47%  "

$ws.Range("B20").Value = "Uniswap"
$ws.Range("C20").Value = "https://coinranking.com/coin/_H5FVG9iW+uniswap-uni"
$ws.Range("D20").Value2 = "'8.00"
$ws.Range("D20").Style = "Normal"
$ws.Range("E20").Value = "  +0.78%  "

$ws.Range("D21").Value2 = "'365.36"
$ws.Range("D21").Style = "Normal"
$ws.Range("E21").Value = "  +2.50%  "

$ws.Range("E22").Value = "  -2.34%  "

$ws.Range("D23").Value2 = "'4.57"
$ws.Range("D23").Style = "Normal"
$ws.Range("E23").Value = "  -2.59%  "

$ws.Range("B24").Value = "Litecoin"
$ws.Range("C24").Value = "https://coinranking.com/coin/D7B1x_ks7WhV5+litecoin-ltc"
$ws.Range("D24").Value2 = "'71.36"
$ws.Range("D24").Style = "Normal"
$ws.Range("E24").Value = "  +1.57%  "

$ws.Range("B25").Value = "Dai"
$ws.Range("C25").Value = "https://coinranking.com/coin/MoTuySvg7+dai-dai"
$ws.Range("D25").Value2 = "'1.00"
$ws.Range("D25").Style = "Normal"
$ws.Range("E25").Value = "  -0.04%  "

$ws.Range("D26").Value2 = "'1.91"
$ws.Range("D26").Style = "Normal"
$ws.Range("E26").Value = "  -6.42%  "

$ws.Range("D27").Value2 = "'9.87"
$ws.Range("D27").Style = "Normal"
$ws.Range("E27").Value = "  -2.18%  "

$ws.Range("D28").Value2 = "'0.998"
$ws.Range("D28").Style = "Normal"
$ws.Range("E28").Value = "  -0.17%  "

$ws.Range("D29").Value = "2.610.38"
$ws.Range("E29").Value = "  -2.85%  "

$ws.Range("D30").Value = "0.0₃0967"
$ws.Range("E30").Value = "  -3.45%  "

$ws.Range("D31").Value2 = "'535.30"
$ws.Range("D31").Style = "Normal"
$ws.Range("E31").Value = "  -0.01%  "

$ws.Range("D32").Value2 = "'8.26"
$ws.Range("D32").Style = "Normal"
$ws.Range("E32").Value = "  -0.31%  "

$ws.Range("E33").Value = "  +0.24%  "

$ws.Range("D34").Value2 = "'1.31"
$ws.Range("D34").Style = "Normal"
$ws.Range("E34").Value = "  -4.47%  "

$ws.Range("D35").Value2 = "'0.999"
$ws.Range("D35").Style = "Normal"
$ws.Range("E35").Value = "  +0.00%  "

$ws.Range("D36").Value2 = "'0.128"
$ws.Range("D36").Style = "Normal"
$ws.Range("E36").Value = "  -4.79%  "

$ws.Range("D37").Value2 = "'158.84"
$ws.Range("D37").Style = "Normal"
$ws.Range("E37").Value = "  +0.95%  "

$ws.Range("E38").Value = "  -3.71%  "

$ws.Range("D39").Value2 = "'18.62"
$ws.Range("D39").Style = "Normal"
$ws.Range("E39").Value = "  -1.14%  "

$ws.Range("D41").Value2 = "'1.79"
$ws.Range("D41").Style = "Normal"
$ws.Range("E41").Value = "  -1.39%  "

$ws.Range("E42").Value = "  -2.41%  "

$ws.Range("D43").Value2 = "'5.13"
$ws.Range("D43").Style = "Normal"
$ws.Range("E43").Value = "  -1.82%  "

$ws.Range("D44").Value2 = "'1.00"
$ws.Range("D44").Style = "Normal"
$ws.Range("E44").Value = "  -0.09%  "

$ws.Range("D45").Value2 = "'2.51"
$ws.Range("D45").Style = "Normal"
$ws.Range("E45").Value = "  -1.98%  "

$ws.Range("D46").Value2 = "'144.96"
$ws.Range("D46").Style = "Normal"
$ws.Range("E46").Value = "  -4.26%  "

$ws.Range("E47").Value = "  -1.43%  "

$ws.Range("D48").Value2 = "'0.549"
$ws.Range("D48").Style = "Normal"
$ws.Range("E48").Value = "  -3.19%  "

$ws.Range("D49").Value = "0.0₆0273"
$ws.Range("E49").Value = "  -3.00%  "

$ws.Range("E50").Value = "  -1.73%  "

$ws.Range("E51").Value = "  -1.95%  "
